$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Station lookup table (Floyd-Warshall path names) written to AE:AF, rows 2-29
# ---------------------------------------------------------------------------
$stations = @(
  "San Isidro",
  "Belenes",
  "Tabachines",
  "Periferico Norte",
  "Independencia Norte",
  "Zoquipan",
  "Plaza Patria",
  "Division del Norte",
  "Vallarta",
  "Colon",
  "Mezquitan",
  "Facultad de Medicina",
  "Seminario",
  "Americas",
  "Juarez",
  "San Juan de Dios",
  "Washington",
  "Agua Azul",
  "Parque Metropolitano",
  "Estampida",
  "Monumento",
  "CUCEI",
  "Unidad Deportiva",
  "El Dean",
  "Abastos",
  "Fray Angelico",
  "Periferico Sur",
  "Adolf Horn"
)

for ($i = 0; $i -lt $stations.Length; $i++) {
  $row = 2 + $i
  $ws.Range("AE$row").Value = $i
  $ws.Range("AE$row").Style = "Normal"
  $ws.Range("AF$row").Value = $stations[$i]
  $ws.Range("AF$row").Style = "Normal"
}

# ---------------------------------------------------------------------------
# New helper-column widths (AD / AE / AF)
# ---------------------------------------------------------------------------
$ws.Columns("AD").ColumnWidth = 2.8333333333333335
$ws.Columns("AE").ColumnWidth = 3.8333333333333335
$ws.Columns("AF").ColumnWidth = 22.666666666666668

# ---------------------------------------------------------------------------
# I30: drop the fill, keep the centered alignment
# ---------------------------------------------------------------------------
$ws.Range("I30").Interior.Pattern = -4142
$ws.Range("I30").HorizontalAlignment = -4108
$ws.Range("I30").VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# Two blank spacer rows below the table (31-32), columns A-AB
# ---------------------------------------------------------------------------
$ws.Range("A31:AB32").Value = 0
$ws.Range("A31:AB32").ClearContents()
$ws.Range("A31:AB32").Style = "Normal"

# ---------------------------------------------------------------------------
# Freeze the header row, restore the active selection
# ---------------------------------------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("T11").Select()

# ---------------------------------------------------------------------------
# Print setup (paper size / orientation)
# ---------------------------------------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
